# Applies the edit described in the commit diff:
#  - Clear the header text in B1 (was "White"), leaving the cell blank
#    but keeping its existing style.
#  - Delete the entire "Total" column I (the first Total column), which
#    shifts the second group of columns (J:O -> I:N) one column left.
#  - Delete the entire "Total" column that is now O (was P, the second
#    Total column), shrinking the used range to A1:N16.
#  - Leave the selection on column O, matching the state left behind
#    after the last column deletion in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Blank out the "White" header text in B1 (style stays, value goes away).
$ws.Range("B1").ClearContents() | Out-Null

# Remove the first "Total" column (I). Everything to its right shifts left.
$ws.Columns("I:I").Delete() | Out-Null

# Remove the second "Total" column, which is now column O after the shift.
$ws.Columns("O:O").Delete() | Out-Null

# Match the final selection state (column O, now past the data).
$ws.Columns("O:O").Select() | Out-Null
